# This ObjTables-generated workbook stamps every sheet's first one or two
# cells with a header string of the form:
#   !!!ObjTables ObjTablesVersion='0.0.8' Date='YYYY-MM-DD HH:MM:SS'
#   !!ObjTables Type='Data' Id='<Name>' Name='<Name>' Date='YYYY-MM-DD HH:MM:SS' ObjTablesVersion='0.0.8'
# The commit re-stamped the whole file with a single new generation
# timestamp. Walk every worksheet, find any of these header cells in A1/A2,
# and rewrite just the embedded Date='...' value, leaving everything else in
# the string untouched.

$wb = $excel.ActiveWorkbook

$newDate = "2020-03-05 18:11:01"
$datePattern = "Date='\d{4}-\d{2}-\d{2} \d{2}:\d{2}:\d{2}'"
$replacement = "Date='$newDate'"

$sheetCount = $wb.Worksheets.Count

for ($i = 1; $i -le $sheetCount; $i++) {
  $ws = $wb.Worksheets.Item($i)

  # The sheets are protected, which blocks direct cell writes - lift
  # protection before editing, then restore it afterwards.
  $ws.Unprotect()

  foreach ($addr in @("A1", "A2")) {
    $cell = $ws.Range($addr)
    $val = $cell.Value2
    if (($val -ne $null) -and ($val -is [string]) -and ($val.StartsWith("!!"))) {
      $newVal = [System.Text.RegularExpressions.Regex]::Replace($val, $datePattern, $replacement)
      if ($newVal -ne $val) {
        $cell.Value2 = $newVal
      }
    }
  }

  # Restore sheet protection (contents/objects/scenarios locked, matching
  # the workbook's original protection intent).
  $ws.Protect($null, $true, $true, $true, $false, $false, $false, $false, $false, $false, $false, $false, $false, $false, $false, $false)
}
